$d = $word.ActiveDocument

# The manuscript's Authors file gained a new "Acknowledgements" paragraph
# (moved here from elsewhere for double-blind review anonymity), appended
# right after the existing "Keywords" paragraph, using the same Body Text
# style as the surrounding metadata paragraphs.

$ackText = "Acknowledgements : We warmly thank Behaviour Interactive Inc. for our collaboration which helps in advancing predator-prey science. We thank Julien C" + [char]0x00E9 + "r" + [char]0x00E9 + " who managed the communication between us and Behaviour Interactive Inc., and Marine Dupuy who helped us prepare the data used in this study. We thank Alastair J. Wilson for insightful comments on the results. We also thank all members of Pierre-Olivier Montiglio" + [char]0x2019 + "s laboratory who provided useful feedback on earlier versions of this work. This work was supported by an NSERC postgraduate doctoral scholarship (569716-2022)."

$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)

$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">' + $ackText + '</w:t></w:r></w:p>'

$insertionPoint.InsertXML($fragment)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
Write-Output $d.Paragraphs.Last.Range.Text
